# Weekly fruit/vegetable price update: a new daily record is inserted
# for "Femacal de La Calera - Espinaca" right before the existing row 174,
# pushing the rest of the data block (old rows 174-215) down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 174 - shifts old rows 174..215 down to 175..216
# and extends the used range to A1:R216.
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row 174 with the new price record.
$ws.Cells.Item(174, 1).Value = 3
$ws.Cells.Item(174, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(174, 3).Value = "Coquimbo"
$ws.Cells.Item(174, 4).Value = 44508
$ws.Cells.Item(174, 5).Value = 5
$ws.Cells.Item(174, 6).Value = 100112012
$ws.Cells.Item(174, 7).Value = "Espinaca"
$ws.Cells.Item(174, 8).Value = "Sin especificar"
$ws.Cells.Item(174, 9).Value = "Primera"
$ws.Cells.Item(174, 10).Value = 150
$ws.Cells.Item(174, 11).Value = 2500
$ws.Cells.Item(174, 12).Value = 2700
$ws.Cells.Item(174, 13).Value = 2607
$ws.Cells.Item(174, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(174, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(174, 16).Value = 869
$ws.Cells.Item(174, 17).Value = 3
$ws.Cells.Item(174, 18).Value = "Hortaliza"

# Match the date-formatted style used by the rest of column D.
$ws.Cells.Item(174, 4).NumberFormat = $ws.Cells.Item(175, 4).NumberFormat
